$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column F is "dSF" - update per repulled data
$ws.Range("F2").Value = 3
$ws.Range("F3").Value = -1
$ws.Range("F4").Value = 5
$ws.Range("F5").Value = -2
